$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new "PortfolioTable" worksheet, positioned before "Landing"
#    (Worksheets.Add() with no args inserts before the active sheet,
#    which for a freshly-opened single-sheet workbook is "Landing").
# ------------------------------------------------------------------
$landing = $wb.Worksheets.Item("Landing")
$newSheet = $wb.Worksheets.Add($landing)
$newSheet.Name = "PortfolioTable"
$ws = $newSheet

# ------------------------------------------------------------------
# 2. Populate the data (header row + 5 data rows, A1:D6)
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Portfolio 1"
$ws.Range("B1").Value = "Portfolio 2"
$ws.Range("C1").Value = "Portfolio 3"
$ws.Range("D1").Value = "Portfolio 4"

$ws.Range("A2").Value = " Berkshire Hathaway CL B"
$ws.Range("B2").Value = " Berkshire Hathaway CL B"
$ws.Range("C2").Value = " Charles Schwab Corp."
$ws.Range("D2").Value = " Charles Schwab Corp."

$ws.Range("A3").Value = " CarMax Inc."
$ws.Range("B3").Value = " CarMax Inc."
$ws.Range("C3").Value = " HEICO Corp. CL A"
$ws.Range("D3").Value = " Taiwan Semiconductor S.A."

$ws.Range("A4").Value = " Markel Corp."
$ws.Range("B4").Value = " Markel Corp."
$ws.Range("C4").Value = " Boot Barn Holdings"
$ws.Range("D4").Value = " Fastenal"

$ws.Range("A5").Value = " Alphabet Inc. CL C"
$ws.Range("B5").Value = " Alphabet Inc. CL C"
$ws.Range("C5").Value = " Levi Strauss & Co."
$ws.Range("D5").Value = " Ametek Inc."

$ws.Range("A6").Value = " Ametek Inc."
$ws.Range("B6").Value = " Ametek Inc."
$ws.Range("C6").Value = " Roper Technologies Inc."
$ws.Range("D6").Value = " Ishares U.S. Home Constructi"

# ------------------------------------------------------------------
# 3. Turn A1:D6 into a table (mirrors the Power Query "PortfolioTable_2"
#    load-to-table that was added upstream)
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:D6"), 0, 1)
$lo.Name = "PortfolioTable_2"
$lo.TableStyle = "TableStyleMedium2"

# ------------------------------------------------------------------
# 4. Column widths to match the sizes used for the lookalike table
#    elsewhere in the workbook.
# ------------------------------------------------------------------
$ws.Range("A1:B1").EntireColumn.ColumnWidth = 22.6
$ws.Range("C1").EntireColumn.ColumnWidth = 22
$ws.Range("D1").EntireColumn.ColumnWidth = 26.4

# ------------------------------------------------------------------
# 5. Defined name ExternalData_1, scoped to the new sheet, hidden --
#    mirrors the query's "load to" external-data range name.
# ------------------------------------------------------------------
$extName = $ws.Names.Add("ExternalData_1", "=PortfolioTable!`$A`$1:`$D`$6")
$extName.Visible = $false

# ------------------------------------------------------------------
# 6. Restore the original ExternalData_2 defined name's sheet-qualified
#    RefersTo text (inserting a sheet shifts localSheetId automatically,
#    just re-assert the formula text so the sheet prefix survives).
# ------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "ExternalData_2") {
        $n.RefersTo = "=Landing!#REF!"
    }
}

# ------------------------------------------------------------------
# 7. Selections / active sheet: Landing no longer the active tab,
#    selection moves to C4; PortfolioTable becomes the active tab with
#    D20 selected.
# ------------------------------------------------------------------
$landing.Activate()
$landing.Range("C4").Select()
$ws.Activate()
$ws.Range("D20").Select()
